$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'Multi-Utilities(18)', 0.6011530979219659),
    @(3, 'Road & Rail(22)', 0.555478110834445),
    @(4, 'Air Freight & Logistics(11)', 0.5120315985335757),
    @(5, 'Electric Utilities(28)', 0.5111832501130005),
    @(6, 'Banks(251)', 0.4661726939731133),
    @(7, 'Building Products(24)', 0.4224192450983748),
    @(8, 'Machinery(86)', 0.4185327478649988),
    @(9, 'Marine(15)', 0.3878966859677347),
    @(10, 'ETF(303)', 0.3793054930628602),
    @(11, 'Gas Utilities(12)', 0.3704453840804315),
    @(12, 'Trading Companies & Distributors(25)', 0.3654504162682767),
    @(13, 'Metals & Mining(106)', 0.3584045941502245),
    @(14, 'Thrifts & Mortgage Finance(47)', 0.3547728109197359),
    @(15, 'Energy Equipment & Services(38)', 0.3434119444041198),
    @(16, 'Equity Real Estate Investment Trusts ...(98)', 0.3194618528503097),
    @(17, 'Insurance(75)', 0.3180148923140174),
    @(18, 'Water Utilities(13)', 0.3159640804598254),
    @(19, 'Auto Components(21)', 0.3055511741165613),
    @(20, 'Chemicals(52)', 0.2955594717296125),
    @(21, 'Construction & Engineering(21)', 0.2795584993604937),
    @(22, 'Life Sciences Tools & Services(19)', 0.2631741957874983),
    @(23, 'Specialty Retail(59)', 0.2332080711894651),
    @(24, 'Capital Markets(76)', 0.2288803717849751),
    @(25, 'Semiconductors & Semiconductor Equipment(70)', 0.2109124772980088),
    @(26, 'Electrical Equipment(28)', 0.1998757633264295),
    @(27, 'Commercial Services & Supplies(52)', 0.1914991355344332),
    @(28, 'Professional Services(35)', 0.1896722757203425),
    @(29, 'Electronic Equipment, Instruments & C...(78)', 0.1835585479573281),
    @(30, 'Aerospace & Defense(37)', 0.1828968893596951),
    @(31, 'Hotels, Restaurants & Leisure(51)', 0.1816312826567509),
    @(32, 'Oil, Gas & Consumable Fuels(125)', 0.178122319223604),
    @(33, 'Health Care Providers & Services(47)', 0.1614771819879285),
    @(34, 'Pharmaceuticals(53)', 0.1528481640262463),
    @(35, 'Communications Equipment(45)', 0.1519806083713743),
    @(36, 'Media(42)', 0.1391667115132491),
    @(37, 'Household Durables(39)', 0.1349024315186565),
    @(38, 'Health Care Equipment & Supplies(86)', 0.1203728078300313),
    @(39, 'Biotechnology(128)', 0.1111411566930368),
    @(40, 'IT Services(52)', 0.09920933229872297),
    @(41, 'Food Products(46)', 0.08887282488690043),
    @(42, 'Software(70)', 0.07670625032558819),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}
